$d = $word.ActiveDocument

# Find the paragraph that reads "Item based:" - the new "Cosine Similarity:"
# section is inserted right before it, reusing the blank paragraph that
# currently sits immediately above it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "Item based:") {
        $target = $d.Paragraphs.Item($i - 1)
        break
    }
}

# Turn the existing blank paragraph into the "Cosine Similarity:" heading.
$target.Range.Text = "Cosine Similarity:"

# Insert a new paragraph right after it holding the explanatory text.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "In the case of cosine similarity, two items are thought of as to be vectors in a k dimensional space. The cosine of the angle between the two vectors is seen to be the similarity between the two items, being a value between 0 and 1. The cosine_sim.py program included in the files is a basic representation on how we computed similarity scores between users and items."
